$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.292.69'
$ws.Range("E2").Value = '  +1.47%  '
$ws.Range("D3").Value = '2.673.71'
$ws.Range("E3").Value = '  +1.56%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '599.62'
$ws.Range("E5").Value = '  -1.17%  '
$ws.Range("D6").Value = '176.59'
$ws.Range("E6").Value = '  -2.57%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.523'
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").Value = '2.674.61'
$ws.Range("E9").Value = '  +1.64%  '
$ws.Range("E10").Value = '  +2.82%  '
$ws.Range("E11").Value = '  +2.09%  '
$ws.Range("D12").Value = '0.353'
$ws.Range("E12").Value = '  +1.96%  '
$ws.Range("D13").Value = '5.02'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("D14").Value = '3.167.09'
$ws.Range("E14").Value = '  +1.70%  '
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '72.161.99'
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").Value = '26.32'
$ws.Range("E17").Value = '  -1.48%  '
$ws.Range("D18").Value = '2.689.86'
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("D19").Value = '12.04'
$ws.Range("E19").Value = '  +4.20%  '
$ws.Range("D20").Value = '8.03'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").Value = '369.84'
$ws.Range("E21").Value = '  -3.12%  '
$ws.Range("D22").Value = '4.16'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = '2.04'
$ws.Range("E23").Value = '  +6.03%  '
$ws.Range("D24").Value = '71.96'
$ws.Range("E24").Value = '  -0.63%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("D26").Value = '4.34'
$ws.Range("E26").Value = '  -2.90%  '
$ws.Range("D27").Value = '9.81'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '2.827.27'
$ws.Range("E28").Value = '  +2.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.17%  '
$ws.Range("E30").Value = '  -2.64%  '
$ws.Range("D31").Value = '8.07'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '510.44'
$ws.Range("E32").Value = '  -5.90%  '
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = '163.54'
$ws.Range("E36").Value = '  -1.50%  '
$ws.Range("D37").Value = '19.55'
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").Value = '19.09'
$ws.Range("E38").Value = '  +0.31%  '
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("E40").Value = '  -3.47%  '
$ws.Range("D41").Value = '0.107'
$ws.Range("E41").Value = '  -9.50%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = '5.01'
$ws.Range("E43").Value = '  -0.96%  '
$ws.Range("E44").Value = '  -2.85%  '
$ws.Range("D45").Value = '0.332'
$ws.Range("E45").Value = '  +0.29%  '
$ws.Range("D46").Value = '39.28'
$ws.Range("E46").Value = '  -1.65%  '
$ws.Range("D47").Value = '153.17'
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").Value = '3.72'
$ws.Range("E48").Value = '  +2.23%  '
$ws.Range("D49").Value = '0.552'
$ws.Range("E49").Value = '  +3.07%  '
$ws.Range("D50").Value = '1.72'
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = '0.0768'
$ws.Range("E51").Value = '  +1.01%  '
